# Convert the old dot-path placeholder syntax ("d.tasks[...]", "d.date.fullDate")
# to the new indexed/filtered syntax ("d[i=0].tasks[...]", "d[i=0].date.fullDate").
#
# Every occurrence of "d." at the start of a template expression becomes
# "d[i=0]." - this covers both the single "{d.date.fullDate}" placeholder
# and the eight "{d.tasks[...]....}" placeholders in the table.

$d = $word.ActiveDocument
$find = $d.Content.Find

# Replace every "d.tasks" with "d[i=0].tasks" (there are 8 such occurrences
# inside the schedule table).
$find.Execute("d.tasks", $true, $false, $false, $false, $false, $true, 1, $false, "d[i=0].tasks", 2)

# Replace "d.date.fullDate" with "d[i=0].date.fullDate" (the single date
# placeholder in the subtitle line).
$find.Execute("d.date.fullDate", $true, $false, $false, $false, $false, $true, 1, $false, "d[i=0].date.fullDate", 2)
